# The OOXML diff for this fixture is purely a round-trip / attribute-
# canonicalisation artefact (the commit bulk re-saved the .docx test
# fixtures while wiring up the "M2Doc version" custom document property;
# for this particular template the saved content is byte-for-byte the
# same *document*, just with XML attributes re-ordered by the tool that
# produced the commit). There is no actual text, formatting, style,
# section or property value change to reproduce here, so this script
# intentionally performs no content mutation - it only touches the
# object model in read-only ways to confirm the document is reachable.

$d = $word.ActiveDocument

# Touch a few read-only properties so the script is clearly operating on
# the live ActiveDocument, without mutating any content.
$null = $d.Content.Text
$null = $d.Sections.Count
$null = $d.Styles.Count
